# Updated symbol list on Tue Feb 14 11:20:35 UTC 2023 with GitHub Actions
#
# Refreshes the cryptocurrency price/volume snapshot in the "Price" (D)
# and "Volume(1h)" (E) columns. Values are written as literal text (a
# leading apostrophe forces text entry) so formats like trailing zeros
# ("5.000") and percent strings ("-0.79%") are preserved exactly, matching
# how the source data was originally stored as inline strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''292.78'
$ws.Range("E2").Value = '''-0.79%'
$ws.Range("D3").Value = '''40.32'
$ws.Range("E3").Value = '''0.61%'
$ws.Range("D4").Value = '''5.000'
$ws.Range("D5").Value = '''0.07327'
$ws.Range("E5").Value = '''-0.77%'
$ws.Range("D6").Value = '''1.556'
$ws.Range("E6").Value = '''1.54%'
$ws.Range("D7").Value = '''0.9221'
$ws.Range("E7").Value = '''0.09%'
$ws.Range("D8").Value = '''2.381'
$ws.Range("E8").Value = '''-0.76%'
$ws.Range("D9").Value = '''0.1187'
$ws.Range("E9").Value = '''1.43%'
$ws.Range("D10").Value = '''0.1815'
$ws.Range("E10").Value = '''3.35%'
$ws.Range("D11").Value = '''0.04394'
$ws.Range("E11").Value = '''5.43%'
$ws.Range("D12").Value = '''0.08797'
$ws.Range("E12").Value = '''2.09%'
$ws.Range("E13").Value = '''-0.05%'
$ws.Range("D14").Value = '''0.001277'
$ws.Range("E14").Value = '''0.46%'
$ws.Range("D15").Value = '''0.005770'
$ws.Range("E15").Value = '''0.08%'
$ws.Range("D16").Value = '''3.342'
$ws.Range("E16").Value = '''-0.95%'
$ws.Range("D17").Value = '''4.292'
$ws.Range("E17").Value = '''-0.45%'
$ws.Range("D18").Value = '''0.3328'
$ws.Range("E18").Value = '''0.93%'
$ws.Range("D19").Value = '''7.882'
$ws.Range("E19").Value = '''3.89%'
$ws.Range("E20").Value = '''2.60%'
$ws.Range("D21").Value = '''0.2802'
$ws.Range("E21").Value = '''-0.31%'
$ws.Range("D22").Value = '''0.03916'
$ws.Range("E22").Value = '''2.30%'
$ws.Range("D23").Value = '''0.001261'
$ws.Range("E23").Value = '''-1.85%'
$ws.Range("D24").Value = '''0.003810'
$ws.Range("E24").Value = '''-2.25%'
$ws.Range("E25").Value = '''-7.35%'
$ws.Range("D26").Value = '''0.0003726'
$ws.Range("E26").Value = '''-0.25%'
$ws.Range("D38").Value = '''0.02336'
$ws.Range("E38").Value = '''1.45%'
$ws.Range("D39").Value = '''0.05075'
$ws.Range("E39").Value = '''1.47%'
$ws.Range("D40").Value = '''0.006017'
$ws.Range("E40").Value = '''49.82%'
$ws.Range("D41").Value = '''0.007810'
$ws.Range("E41").Value = '''1.31%'
$ws.Range("E42").Value = '''1.29%'
$ws.Range("D43").Value = '''0.007388'
$ws.Range("E43").Value = '''-0.42%'
$ws.Range("D44").Value = '''0.008050'
$ws.Range("E44").Value = '''4.72%'
$ws.Range("D45").Value = '''0.2907'
$ws.Range("E45").Value = '''-9.35%'
$ws.Range("D46").Value = '''0.00006211'
$ws.Range("E46").Value = '''-3.92%'
$ws.Range("D47").Value = '''0.00000000751'
$ws.Range("E47").Value = '''-0.24%'
$ws.Range("D49").Value = '''0.004205'
$ws.Range("E49").Value = '''-0.24%'
$ws.Range("D50").Value = '''0.00002102'
$ws.Range("E50").Value = '''-0.24%'
$ws.Range("D51").Value = '''0.0002002'
$ws.Range("E51").Value = '''-0.24%'
